$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column letters in the order the values are stored in each row array below:
# E, G, H, I, J, K, M, N, O, P, Q, R, S, T
$cols = @(5,7,8,9,10,11,13,14,15,16,17,18,19,20)

$rows = @(
    @(3.0,6.423576,19.270728,0.001681024218962088,0.001681024218962088,3.0,239.0839323333333,717.251797,0.4086975387666237,0.4086975387666237,1535.773809722024,13821.96428749821,0.0006870304608968913,0.0006870304608968913),
    @(3.0,6.423576,19.270728,0.001681024218962088,0.001681024218962088,3.0,117.0512696666667,351.153809,0.2000910950200451,0.2000910950200451,751.8877266003279,6766.989539402952,0.0003363579767273402,0.0003363579767273402),
    @(3.0,6.423576,19.270728,0.001681024218962088,0.001681024218962088,3.0,171.15883,513.47649,0.2925842480357353,0.2925842480357353,1099.45175257608,9895.065773184719,0.0004918412070348818,0.0004918412070348818),
    @(3.0,6.423576,19.270728,0.001681024218962088,0.001681024218962088,3.0,57.695868,173.087604,0.09862711817759588,0.09862711817759588,370.613792983968,3335.524136855712,0.0001657945743029747,0.0001657945743029747),
    @(3.0,3580.644531333333,10741.933594,0.9370403925578976,0.9370403925578976,3.0,239.0839323333333,717.251797,0.4086975387666237,0.4086975387666237,856074.5748390188,7704671.173551168,0.3829661021633236,0.3829661021633236),
    @(3.0,3580.644531333333,10741.933594,0.9370403925578976,0.9370403925578976,3.0,117.0512696666667,351.153809,0.2000910950200451,0.2000910950200451,419118.9886175733,3772070.89755816,0.1874934382249226,0.1874934382249226),
    @(3.0,3580.644531333333,10741.933594,0.9370403925578976,0.9370403925578976,3.0,171.15883,513.47649,0.2925842480357353,0.2925842480357353,612858.9286289116,5515730.357660205,0.2741632586356627,0.2741632586356627),
    @(3.0,3580.644531333333,10741.933594,0.9370403925578976,0.9370403925578976,3.0,57.695868,173.087604,0.09862711817759588,0.09862711817759588,206588.3942347298,1859295.548112569,0.09241759353398861,0.09241759353398861),
    @(3.0,0.9157713333333334,2.747314,0.0002396537054071653,0.0002396537054071653,3.0,239.0839323333333,717.251797,0.4086975387666237,0.4086975387666237,218.9462114914731,1970.515903423258,0.00009794587955620994,0.00009794587955620994),
    @(3.0,0.9157713333333334,2.747314,0.0002396537054071653,0.0002396537054071653,3.0,117.0512696666667,351.153809,0.2000910950200451,0.2000910950200451,107.1921972910029,964.7297756190261,0.00004795257234053099,0.000047952572340531),
    @(3.0,0.9157713333333334,2.747314,0.0002396537054071653,0.0002396537054071653,3.0,171.15883,513.47649,0.2925842480357353,0.2925842480357353,156.7423499608733,1410.68114964786,0.00007011889918553307,0.00007011889918553309),
    @(3.0,0.9157713333333334,2.747314,0.0002396537054071653,0.0002396537054071653,3.0,57.695868,173.087604,0.09862711817759588,0.09862711817759588,52.836221966184,475.525997695656,0.00002363635432489124,0.00002363635432489124),
    @(3.0,233.243637,699.7309110000001,0.0610389295177331,0.06103892951773311,3.0,239.0839323333333,717.251797,0.4086975387666237,0.4086975387666237,55764.80592568858,501883.2533311971,0.02494646026284693,0.02494646026284694),
    @(3.0,233.243637,699.7309110000001,0.0610389295177331,0.06103892951773311,3.0,117.0512696666667,351.153809,0.2000910950200451,0.2000910950200451,27301.46385252112,245713.17467269,0.01221334624605457,0.01221334624605457),
    @(3.0,233.243637,699.7309110000001,0.0610389295177331,0.06103892951773311,3.0,171.15883,513.47649,0.2925842480357353,0.2925842480357353,39921.70801386471,359295.3721247825,0.01785902929385218,0.01785902929385219),
    @(3.0,233.243637,699.7309110000001,0.0610389295177331,0.06103892951773311,3.0,57.695868,173.087604,0.09862711817759588,0.09862711817759588,13457.19409219192,121114.7468297273,0.006020093714979408,0.006020093714979409)
)

$startRow = 2
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $rowVals = $rows[$i]
    for ($j = 0; $j -lt $cols.Count; $j++) {
        $c = $cols[$j]
        $ws.Cells.Item($r, $c).Value = $rowVals[$j]
    }
}
